# Update CRP/"profit" summary columns (H-N) across ALC/CRP/CUL/GSM/LTW/WVR sheets
# per scheduled Ixion_Profits recalculation.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 62673
$ws.Range("I21").Value = 59009.5
$ws.Range("J21").Value = 70000
$ws.Range("K21").Value = 59009.5
$ws.Range("L21").Value = 70000
$ws.Range("M21").Value = -58541.5
$ws.Range("N21").Value = -70936

$ws.Range("H23").Value = 62673
$ws.Range("I23").Value = 59009.5
$ws.Range("J23").Value = 70000
$ws.Range("K23").Value = 59009.5
$ws.Range("L23").Value = 70000
$ws.Range("M23").Value = -58775.5
$ws.Range("N23").Value = -70468

$ws.Range("H49").Value = 4017
$ws.Range("I49").Value = 4017
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 12051
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -11915
$ws.Range("N49").ClearContents()

$ws.Range("H76").Value = 8580.6
$ws.Range("I76").Value = 10647.077
$ws.Range("J76").Value = 4742.857
$ws.Range("K76").Value = 10647.077
$ws.Range("L76").Value = 4742.857
$ws.Range("M76").Value = -10332.077
$ws.Range("N76").Value = -5372.857

$ws.Range("H79").Value = 8580.6
$ws.Range("I79").Value = 10647.077
$ws.Range("J79").Value = 4742.857
$ws.Range("K79").Value = 10647.077
$ws.Range("L79").Value = 4742.857
$ws.Range("M79").Value = -9555.076999999999
$ws.Range("N79").Value = -6926.857

$ws.Range("H127").Value = 1843.8064
$ws.Range("I127").Value = 695.5
$ws.Range("J127").Value = 2119.4
$ws.Range("K127").Value = 2086.5
$ws.Range("L127").Value = 6358.200000000001
$ws.Range("M127").Value = 2873.5
$ws.Range("N127").Value = -16278.2

$ws.Range("H138").Value = 5055.373
$ws.Range("J138").Value = 8487.468999999999
$ws.Range("L138").Value = 25462.407
$ws.Range("N138").Value = -35742.407

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 3500
$ws.Range("I35").Value = 3500
$ws.Range("K35").Value = 3500
$ws.Range("M35").Value = -3206

$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H74").Value = 20966
$ws.Range("J74").Value = 20966
$ws.Range("L74").Value = 20966
$ws.Range("N74").Value = -22714

$ws.Range("H77").Value = 20966
$ws.Range("J77").Value = 20966
$ws.Range("L77").Value = 62898
$ws.Range("N77").Value = -71634

$ws.Range("H94").Value = 3951.25
$ws.Range("I94").Value = 5151
$ws.Range("J94").Value = 3437.0715
$ws.Range("K94").Value = 5151
$ws.Range("L94").Value = 3437.0715
$ws.Range("M94").Value = -4700
$ws.Range("N94").Value = -4339.0715

$ws.Range("H99").Value = 13894078
$ws.Range("I99").Value = 2579.8
$ws.Range("K99").Value = 2579.8
$ws.Range("M99").Value = -1081.8

$ws.Range("H126").Value = 13894078
$ws.Range("I126").Value = 2579.8
$ws.Range("K126").Value = 7739.400000000001
$ws.Range("M126").Value = -5269.400000000001

$ws.Range("H132").Value = 2119
$ws.Range("I132").Value = 1892.081
$ws.Range("K132").Value = 5676.242999999999
$ws.Range("M132").Value = -3146.242999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2690.4546
$ws.Range("J34").Value = 3155
$ws.Range("L34").Value = 9465
$ws.Range("N34").Value = -9633

$ws.Range("H64").Value = 3400
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3400
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 10200
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -10740

$ws.Range("H67").Value = 3400
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3400
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 10200
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -12072

$ws.Range("H92").Value = 749.3043
$ws.Range("I92").Value = 515
$ws.Range("J92").Value = 798.6316
$ws.Range("K92").Value = 1545
$ws.Range("L92").Value = 2395.8948
$ws.Range("M92").Value = -297
$ws.Range("N92").Value = -4891.8948

$ws.Range("H129").Value = 1893.35
$ws.Range("I129").Value = 1238.375
$ws.Range("J129").Value = 2330
$ws.Range("K129").Value = 3715.125
$ws.Range("L129").Value = 6990
$ws.Range("M129").Value = 1284.875
$ws.Range("N129").Value = -16990

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 3000
$ws.Range("I31").Value = 3000
$ws.Range("K31").Value = 3000
$ws.Range("M31").Value = -2708

$ws.Range("H37").Value = 3000
$ws.Range("I37").Value = 3000
$ws.Range("K37").Value = 3000
$ws.Range("M37").Value = -2723

$ws.Range("H62").Value = 22500
$ws.Range("J62").Value = 22500
$ws.Range("L62").Value = 22500
$ws.Range("N62").Value = -23872

$ws.Range("H65").Value = 22500
$ws.Range("J65").Value = 22500
$ws.Range("L65").Value = 67500
$ws.Range("N65").Value = -74364

$ws.Range("H70").Value = 5772.528
$ws.Range("I70").Value = 5973.154
$ws.Range("J70").Value = 5250.9
$ws.Range("K70").Value = 5973.154
$ws.Range("L70").Value = 5250.9
$ws.Range("M70").Value = -5703.154
$ws.Range("N70").Value = -5790.9

$ws.Range("H73").Value = 5772.528
$ws.Range("I73").Value = 5973.154
$ws.Range("J73").Value = 5250.9
$ws.Range("K73").Value = 5973.154
$ws.Range("L73").Value = 5250.9
$ws.Range("M73").Value = -5037.154
$ws.Range("N73").Value = -7122.9

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 22161.666
$ws.Range("J63").Value = 22161.666
$ws.Range("L63").Value = 22161.666
$ws.Range("N63").Value = -23659.666

$ws.Range("H66").Value = 22161.666
$ws.Range("J66").Value = 22161.666
$ws.Range("L66").Value = 66484.99800000001
$ws.Range("N66").Value = -73972.99800000001

$ws.Range("H132").Value = 9851186
$ws.Range("I132").Value = 12747688
$ws.Range("J132").Value = 3079.3
$ws.Range("K132").Value = 38243064
$ws.Range("L132").Value = 9237.900000000001
$ws.Range("M132").Value = -38240534
$ws.Range("N132").Value = -14297.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1796.6
$ws.Range("I132").Value = 1158.1177
$ws.Range("J132").Value = 3153.375
$ws.Range("K132").Value = 3474.3531
$ws.Range("L132").Value = 9460.125
$ws.Range("M132").Value = -944.3531000000003
$ws.Range("N132").Value = -14520.125
